$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3782.5571
$ws.Range("I64").Value = 3610.8108
$ws.Range("J64").Value = 3975.121
$ws.Range("K64").Value = 3610.8108
$ws.Range("L64").Value = 3975.121
$ws.Range("M64").Value = -3362.8108
$ws.Range("N64").Value = -4471.121

$ws.Range("H67").Value = 3782.5571
$ws.Range("I67").Value = 3610.8108
$ws.Range("J67").Value = 3975.121
$ws.Range("K67").Value = 3610.8108
$ws.Range("L67").Value = 3975.121
$ws.Range("M67").Value = -2752.8108
$ws.Range("N67").Value = -5691.121

$ws.Range("H76").Value = 3241.8928
$ws.Range("I76").Value = 2807.7273
$ws.Range("J76").Value = 3522.8235
$ws.Range("K76").Value = 2807.7273
$ws.Range("L76").Value = 3522.8235
$ws.Range("M76").Value = -2492.7273
$ws.Range("N76").Value = -4152.8235

$ws.Range("H79").Value = 3241.8928
$ws.Range("I79").Value = 2807.7273
$ws.Range("J79").Value = 3522.8235
$ws.Range("K79").Value = 2807.7273
$ws.Range("L79").Value = 3522.8235
$ws.Range("M79").Value = -1715.7273
$ws.Range("N79").Value = -5706.8235

$ws.Range("H103").Value = 536.2857
$ws.Range("I103").Value = 563
$ws.Range("J103").Value = 500.66666
$ws.Range("K103").Value = 1689
$ws.Range("L103").Value = 1501.99998
$ws.Range("M103").Value = -1103
$ws.Range("N103").Value = -2673.99998

$ws.Range("H116").Value = 3346972.8
$ws.Range("I116").Value = 25642526
$ws.Range("J116").Value = 2639.9
$ws.Range("K116").Value = 25642526
$ws.Range("L116").Value = 2639.9
$ws.Range("M116").Value = -25639084
$ws.Range("N116").Value = -9523.9

$ws.Range("H137").Value = 9279497
$ws.Range("I137").Value = 19231634
$ws.Range("J137").Value = 38227.68
$ws.Range("K137").Value = 57694902
$ws.Range("L137").Value = 114683.04
$ws.Range("M137").Value = -57692352
$ws.Range("N137").Value = -119783.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2600
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 2850
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 2850
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -4222

$ws.Range("H66").Value = 2600
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 2850
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 14250
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -21114

$ws.Range("H109").Value = 34994.5
$ws.Range("J109").Value = 34994.5
$ws.Range("L109").Value = 34994.5
$ws.Range("N109").Value = -37768.5

$ws.Range("H112").Value = 9000
$ws.Range("J112").Value = 9000
$ws.Range("L112").Value = 9000
$ws.Range("N112").Value = -11954

$ws.Range("H132").Value = 2690.8572
$ws.Range("I132").Value = 1929.5
$ws.Range("J132").Value = 3706
$ws.Range("K132").Value = 5788.5
$ws.Range("L132").Value = 11118
$ws.Range("M132").Value = -3258.5
$ws.Range("N132").Value = -16178

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1749.6666
$ws.Range("I86").Value = 1539.6
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 1539.6
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -416.5999999999999
$ws.Range("N86").Value = -5046

$ws.Range("H89").Value = 1749.6666
$ws.Range("I89").Value = 1539.6
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 7698
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -2082
$ws.Range("N89").Value = -25232

$ws.Range("H99").Value = 4931.107
$ws.Range("I99").Value = 5739.476
$ws.Range("J99").Value = 2506
$ws.Range("K99").Value = 5739.476
$ws.Range("L99").Value = 2506
$ws.Range("M99").Value = -4241.476
$ws.Range("N99").Value = -5502

$ws.Range("H105").Value = 2322.2
$ws.Range("I105").Value = 1200
$ws.Range("J105").Value = 2446.889
$ws.Range("K105").Value = 1200
$ws.Range("L105").Value = 2446.889
$ws.Range("M105").Value = 547
$ws.Range("N105").Value = -5940.889

$ws.Range("H128").Value = 4470
$ws.Range("I128").Value = 4470
$ws.Range("K128").Value = 13410
$ws.Range("M128").Value = -10920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2684.5818
$ws.Range("I31").Value = 1478.5
$ws.Range("J31").Value = 5380.5293
$ws.Range("K31").Value = 1478.5
$ws.Range("L31").Value = 5380.5293
$ws.Range("M31").Value = -1183.5
$ws.Range("N31").Value = -5970.5293

$ws.Range("H34").Value = 2684.5818
$ws.Range("I34").Value = 1478.5
$ws.Range("J34").Value = 5380.5293
$ws.Range("K34").Value = 1478.5
$ws.Range("L34").Value = 5380.5293
$ws.Range("M34").Value = -1276.5
$ws.Range("N34").Value = -5784.5293

$ws.Range("H62").Value = 2359
$ws.Range("I62").Value = 2375
$ws.Range("J62").Value = 2348.3333
$ws.Range("K62").Value = 2375
$ws.Range("L62").Value = 2348.3333
$ws.Range("M62").Value = -1751
$ws.Range("N62").Value = -3596.3333

$ws.Range("H65").Value = 2359
$ws.Range("I65").Value = 2375
$ws.Range("J65").Value = 2348.3333
$ws.Range("K65").Value = 11875
$ws.Range("L65").Value = 11741.6665
$ws.Range("M65").Value = -8755
$ws.Range("N65").Value = -17981.6665

$ws.Range("H105").Value = 561
$ws.Range("I105").Value = 499.875
$ws.Range("K105").Value = 499.875
$ws.Range("M105").Value = 1247.125

$ws.Range("H134").Value = 2199.7659
$ws.Range("I134").Value = 1284.5358
$ws.Range("J134").Value = 3548.5264
$ws.Range("K134").Value = 3853.6074
$ws.Range("L134").Value = 10645.5792
$ws.Range("M134").Value = -1318.6074
$ws.Range("N134").Value = -15715.5792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 625447
$ws.Range("J107").Value = 1111748.1
$ws.Range("L107").Value = 3335244.3
$ws.Range("N107").Value = -3339084.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12500
$ws.Range("I70").Value = 27000
$ws.Range("K70").Value = 27000
$ws.Range("M70").Value = -26730

$ws.Range("H73").Value = 12500
$ws.Range("I73").Value = 27000
$ws.Range("K73").Value = 27000
$ws.Range("M73").Value = -26064

$ws.Range("H80").Value = 2879.8
$ws.Range("I80").Value = 2727.1428
$ws.Range("J80").Value = 3013.375
$ws.Range("K80").Value = 2727.1428
$ws.Range("L80").Value = 3013.375
$ws.Range("M80").Value = -1729.1428
$ws.Range("N80").Value = -5009.375

$ws.Range("H83").Value = 2879.8
$ws.Range("I83").Value = 2727.1428
$ws.Range("J83").Value = 3013.375
$ws.Range("K83").Value = 13635.714
$ws.Range("L83").Value = 15066.875
$ws.Range("M83").Value = -8643.714
$ws.Range("N83").Value = -25050.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 8041.4287
$ws.Range("J110").Value = 8041.4287
$ws.Range("L110").Value = 8041.4287
$ws.Range("N110").Value = -16221.4287
